$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q1" fund-holdings sheet, positioned right before
#    the "总计" sheet. We clone an existing fund sheet (2021-Q3, which has
#    the same 7-column layout and the closest row count) so the new sheet
#    inherits identical formatting (header style, index-column style, etc).
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q3")
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet.Copy($totalSheet)
$q1 = $wb.ActiveSheet
$q1.Name = "2022-Q1"

# The template sheet only has 6 rows (header + 5 funds); we need 7 rows
# (header + 6 funds), so add row 7 by duplicating the format of row 6.
$q1.Range("A6").Copy()
$q1.Range("A7").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'001302"
$q1.Range("C2").Value = "前海开源金银珠宝主题精选混合A"
$q1.Range("D2").Value = "'8.61"
$q1.Range("E2").Value = "'91.91"
$q1.Range("F2").Value = "'9.14"
$q1.Range("G2").Value = "'0.7870"
$q1.Range("H2").Value = 5

# Row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'003304"
$q1.Range("C3").Value = "前海开源沪港深核心资源灵活配置混合A"
$q1.Range("D3").Value = "'5.91"
$q1.Range("E3").Value = "'93.10"
$q1.Range("F3").Value = "'7.27"
$q1.Range("G3").Value = "'0.4297"
$q1.Range("H3").Value = 9

# Row 4
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'002207"
$q1.Range("C4").Value = "前海开源金银珠宝主题精选混合C"
$q1.Range("D4").Value = "'3.45"
$q1.Range("E4").Value = "'91.91"
$q1.Range("F4").Value = "'9.14"
$q1.Range("G4").Value = "'0.3153"
$q1.Range("H4").Value = 5

# Row 5
$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "'003305"
$q1.Range("C5").Value = "前海开源沪港深核心资源灵活配置混合C"
$q1.Range("D5").Value = "'2.19"
$q1.Range("E5").Value = "'93.10"
$q1.Range("F5").Value = "'7.27"
$q1.Range("G5").Value = "'0.1592"
$q1.Range("H5").Value = 9

# Row 6
$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "'008842"
$q1.Range("C6").Value = "同泰远见灵活配置混合A"
$q1.Range("D6").Value = "'0.43"
$q1.Range("E6").Value = "'94.53"
$q1.Range("F6").Value = "'3.01"
$q1.Range("G6").Value = "'0.0129"
$q1.Range("H6").Value = 10

# Row 7
$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "'008843"
$q1.Range("C7").Value = "同泰远见灵活配置混合C"
$q1.Range("D7").Value = "'0.09"
$q1.Range("E7").Value = "'94.53"
$q1.Range("F7").Value = "'3.01"
$q1.Range("G7").Value = "'0.0027"
$q1.Range("H7").Value = 10

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new first data row for the
#    2022-Q1 quarter, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()

# Restore formatting on the newly inserted row's cells (Insert leaves the
# index cell A2 unstyled and copies stray formatting into B2:D2).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B2:D2").ClearFormats()

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 1.71

# Renumber the index column for the rows that shifted down.
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
